# Trade #5 closed at 2026-02-17 19:44:07 - unknown UNKNOWN +0.000%
# Updates summary/strategy-status aggregates and appends the new trade row
# to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1300.09
$summary.Range("B4").Value = 0.09
$summary.Range("B5").Value = 0.36
$summary.Range("B6").Value = 5
$summary.Range("B7").Value = 3
$summary.Range("B9").Value = 60

# ---------------------------------------------------------------------
# Strategy Status sheet (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.09
$status.Range("D4").Value = 5
$status.Range("E4").Value = 0.09
$status.Range("F4").Value = 0.09
$status.Range("G4").Value = 60

# ---------------------------------------------------------------------
# New trade row (Trade #5) shared by "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 6

    $ws.Cells.Item($row, 1).Value = 5

    # "2026-02-17" looks like a date, so force it to stay plain text
    # (matches the rest of the Date column, which is inline/string data).
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "19:44:01"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.727273
    $ws.Cells.Item($row, 7).Value = 0.78
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 7.25
    $ws.Cells.Item($row, 10).Value = 0.05
    $ws.Cells.Item($row, 11).Value = 100.09
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}
